$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A94").Value = "GRT-USD"
